$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "97.329.62"
$ws.Range("E2").Value = "  -2.02%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.393.57"
$ws.Range("E3").Value = "  +3.18%  "

$ws.Range("E4").Value = "  +0.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "253.94"
$ws.Range("E5").Value = "  -0.21%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "644.77"
$ws.Range("E6").Value = "  +3.25%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.44"
$ws.Range("E7").Value = "  -0.52%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.419"
$ws.Range("E8").Value = "  +4.54%  "

$ws.Range("B9").Value = "USDC"
$ws.Range("C9").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  +0.07%  "

$ws.Range("B10").Value = "Cardano"
$ws.Range("C10").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.04"
$ws.Range("E10").Value = "  +5.75%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.396.41"
$ws.Range("E11").Value = "  +3.36%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.210"
$ws.Range("E12").Value = "  +4.68%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "41.15"
$ws.Range("E13").Value = "  +3.50%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.21"
$ws.Range("E14").Value = "  +13.04%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000255"
$ws.Range("E15").Value = "  +2.46%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "97.037.97"
$ws.Range("E16").Value = "  -1.98%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.023.43"
$ws.Range("E17").Value = "  +2.95%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.41"
$ws.Range("E18").Value = "  +32.55%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.391.34"
$ws.Range("E19").Value = "  +2.95%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.22"
$ws.Range("E20").Value = "  +12.70%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.72"
$ws.Range("E21").Value = "  +15.24%  "

$ws.Range("B22").Value = "Stellar"
$ws.Range("C22").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.477"
$ws.Range("E22").Value = "  +42.56%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.39"
$ws.Range("E23").Value = "  -1.99%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "499.92"
$ws.Range("E24").Value = "  +2.37%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000201"
$ws.Range("E25").Value = "  -0.24%  "

$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "97.52"
$ws.Range("E26").Value = "  +9.34%  "

$ws.Range("B27").Value = "NEARProtocol"
$ws.Range("C27").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.02"
$ws.Range("E27").Value = "  +6.92%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.54"
$ws.Range("E28").Value = "  +3.80%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.581.00"
$ws.Range("E29").Value = "  +3.43%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.151"
$ws.Range("E30").Value = "  +9.57%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.197"
$ws.Range("E31").Value = "  +3.74%  "

$ws.Range("B32").Value = "Dai"
$ws.Range("C32").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.996"
$ws.Range("E32").Value = "  -0.31%  "

$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.17"
$ws.Range("E33").Value = "  +7.55%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.997"
$ws.Range("E34").Value = "  -0.47%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.562"
$ws.Range("E35").Value = "  +17.65%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "29.39"
$ws.Range("E36").Value = "  +4.99%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.22"
$ws.Range("E37").Value = "  +14.50%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.60"
$ws.Range("E38").Value = "  +4.99%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.152"
$ws.Range("E39").Value = "  +1.11%  "

$ws.Range("E40").Value = "  +12.13%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "503.71"
$ws.Range("E41").Value = "  +3.17%  "

$ws.Range("E42").Value = "  -0.28%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.850"
$ws.Range("E43").Value = "  +9.84%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.63"
$ws.Range("E44").Value = "  -3.32%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0407"
$ws.Range("E45").Value = "  +20.58%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.42"
$ws.Range("E46").Value = "  +13.78%  "

$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.00"
$ws.Range("E47").Value = "  +0.02%  "

$ws.Range("B48").Value = "dogwifhat"
$ws.Range("C48").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.19"
$ws.Range("E48").Value = "  +2.71%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.10"
$ws.Range("E49").Value = "  +11.02%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.55"
$ws.Range("E50").Value = "  +13.53%  "

$ws.Range("B51").Value = "OKB"
$ws.Range("C51").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.40"
$ws.Range("E51").Value = "  +10.93%  "
